$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 17949
$ws.Range("E2").Value = 823
$ws.Range("F2").Value = 823
$ws.Range("G2").Value = 784
$ws.Range("H2").Value = 592
$ws.Range("I2").Value = 592
$ws.Range("J2").Value = -1
$ws.Range("K2").Value = 8971
$ws.Range("L2").Value = 4185
$ws.Range("M2").Value = 4786
$ws.Range("N2").Value = 4784
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 193
$ws.Range("Q2").Value = 1007
$ws.Range("R2").Value = -207
$ws.Range("S2").Value = -809
$ws.Range("T2").Value = 378
$ws.Range("U2").Value = 628
$ws.Range("V2").Value = 1839
$ws.Range("W2").Value = 4.58
$ws.Range("X2").Value = 3.3
$ws.Range("Y2").Value = 13.09
$ws.Range("Z2").Value = 6.53
$ws.Range("AA2").Value = 87.42
$ws.Range("AB2").Value = 2419.7
$ws.Range("AC2").Value = 15350
$ws.Range("AD2").Value = 18.99
$ws.Range("AE2").Value = 123964
$ws.Range("AF2").Value = 2.35
$ws.Range("AG2").Value = 3000
$ws.Range("AH2").Value = 1.03
$ws.Range("AI2").Value = 19.54
$ws.Range("AJ2").Value = 3859124

$ws.Range("D3").Value = 19310
$ws.Range("E3").Value = 771
$ws.Range("F3").Value = 771
$ws.Range("G3").Value = 729
$ws.Range("H3").Value = 568
$ws.Range("I3").Value = 568
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 9524
$ws.Range("L3").Value = 4302
$ws.Range("M3").Value = 5222
$ws.Range("N3").Value = 5221
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 193
$ws.Range("Q3").Value = 890
$ws.Range("R3").Value = -571
$ws.Range("S3").Value = -250
$ws.Range("T3").Value = 397
$ws.Range("U3").Value = 493
$ws.Range("V3").Value = 1858
$ws.Range("W3").Value = 3.99
$ws.Range("X3").Value = 2.94
$ws.Range("Y3").Value = 11.36
$ws.Range("Z3").Value = 6.14
$ws.Range("AA3").Value = 82.39
$ws.Range("AB3").Value = 2654.31
$ws.Range("AC3").Value = 14731
$ws.Range("AD3").Value = 26.68
$ws.Range("AE3").Value = 135283
$ws.Range("AF3").Value = 2.91
$ws.Range("AG3").Value = 2500
$ws.Range("AH3").Value = 0.64
$ws.Range("AI3").Value = 16.97
$ws.Range("AJ3").Value = 3859124

$ws.Range("D4").Value = 22413
$ws.Range("E4").Value = 733
$ws.Range("F4").Value = 733
$ws.Range("G4").Value = 648
$ws.Range("H4").Value = 504
$ws.Range("I4").Value = 503
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 10854
$ws.Range("L4").Value = 5335
$ws.Range("M4").Value = 5519
$ws.Range("N4").Value = 5517
$ws.Range("O4").Value = 2
$ws.Range("P4").Value = 193
$ws.Range("Q4").Value = 1046
$ws.Range("R4").Value = -544
$ws.Range("S4").Value = -174
$ws.Range("T4").Value = 394
$ws.Range("U4").Value = 652
$ws.Range("V4").Value = 1941
$ws.Range("W4").Value = 3.27
$ws.Range("X4").Value = 2.25
$ws.Range("Y4").Value = 9.380000000000001
$ws.Range("Z4").Value = 4.95
$ws.Range("AA4").Value = 96.66
$ws.Range("AB4").Value = 2837.34
$ws.Range("AC4").Value = 13044
$ws.Range("AD4").Value = 15.18
$ws.Range("AE4").Value = 142964
$ws.Range("AF4").Value = 1.38
$ws.Range("AG4").Value = 2500
$ws.Range("AH4").Value = 1.26
$ws.Range("AI4").Value = 19.17
$ws.Range("AJ4").Value = 3859124

$ws.Range("D5").Value = 25526
$ws.Range("E5").Value = 724
$ws.Range("F5").Value = 724
$ws.Range("G5").Value = 673
$ws.Range("H5").Value = 514
$ws.Range("I5").Value = 514
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 13691
$ws.Range("L5").Value = 7759
$ws.Range("M5").Value = 5932
$ws.Range("N5").Value = 5931
$ws.Range("O5").Value = 2
$ws.Range("P5").Value = 193
$ws.Range("Q5").Value = 497
$ws.Range("R5").Value = -2329
$ws.Range("S5").Value = 1801
$ws.Range("T5").Value = 903
$ws.Range("U5").Value = -407
$ws.Range("V5").Value = 4226
$ws.Range("W5").Value = 2.84
$ws.Range("X5").Value = 2.01
$ws.Range("Y5").Value = 8.970000000000001
$ws.Range("Z5").Value = 4.19
$ws.Range("AA5").Value = 130.79
$ws.Range("AB5").Value = 3052.59
$ws.Range("AC5").Value = 13313
$ws.Range("AD5").Value = 17.58
$ws.Range("AE5").Value = 153682
$ws.Range("AF5").Value = 1.52
$ws.Range("AG5").Value = 2500
$ws.Range("AH5").Value = 1.07
$ws.Range("AI5").Value = 18.78
$ws.Range("AJ5").Value = 3859124

$ws.Range("D6").Value = 28025
$ws.Range("E6").Value = 872
$ws.Range("F6").Value = 872
$ws.Range("G6").Value = 714
$ws.Range("H6").Value = 567
$ws.Range("I6").Value = 567
$ws.Range("K6").Value = 14036
$ws.Range("L6").Value = 7711
$ws.Range("M6").Value = 6325
$ws.Range("N6").Value = 6325
$ws.Range("P6").Value = 193
$ws.Range("Q6").Value = 261
$ws.Range("R6").Value = -843
$ws.Range("S6").Value = 164
$ws.Range("T6").Value = 897
$ws.Range("U6").Value = -636
$ws.Range("V6").Value = 4630
$ws.Range("W6").Value = 3.11
$ws.Range("X6").Value = 2.02
$ws.Range("Y6").Value = 9.26
$ws.Range("Z6").Value = 4.09
$ws.Range("AA6").Value = 121.92
$ws.Range("AB6").Value = 3298.77
$ws.Range("AC6").Value = 14699
$ws.Range("AD6").Value = 19.56
$ws.Range("AE6").Value = 163896
$ws.Range("AF6").Value = 1.75
$ws.Range("AG6").Value = 3000
$ws.Range("AH6").Value = 1.04
$ws.Range("AI6").Value = 20.41
$ws.Range("AJ6").Value = 3859124

$ws.Range("D7").Value = 30019
$ws.Range("E7").Value = 984
$ws.Range("G7").Value = 810
$ws.Range("H7").Value = 632
$ws.Range("I7").Value = 645
$ws.Range("K7").Value = 15063
$ws.Range("L7").Value = 8227
$ws.Range("M7").Value = 6836
$ws.Range("N7").Value = 6836
$ws.Range("P7").Value = 191
$ws.Range("Q7").Value = 1185
$ws.Range("R7").Value = -1113
$ws.Range("S7").Value = -161
$ws.Range("T7").Value = 525
$ws.Range("U7").Value = 160
$ws.Range("W7").Value = 3.28
$ws.Range("X7").Value = 2.1
$ws.Range("Y7").Value = 9.800000000000001
$ws.Range("Z7").Value = 4.34
$ws.Range("AA7").Value = 120.34
$ws.Range("AC7").Value = 16703
$ws.Range("AD7").Value = 12.57
$ws.Range("AE7").Value = 177147
$ws.Range("AF7").Value = 1.19
$ws.Range("AG7").Value = 3000
$ws.Range("AH7").Value = 1.43
$ws.Range("AI7").Value = 17.96

$ws.Range("D8").Value = 32198
$ws.Range("E8").Value = 1097
$ws.Range("G8").Value = 947
$ws.Range("H8").Value = 726
$ws.Range("I8").Value = 736
$ws.Range("K8").Value = 15843
$ws.Range("L8").Value = 8394
$ws.Range("M8").Value = 7448
$ws.Range("N8").Value = 7448
$ws.Range("P8").Value = 191
$ws.Range("Q8").Value = 1139
$ws.Range("R8").Value = -784
$ws.Range("S8").Value = -218
$ws.Range("T8").Value = 674
$ws.Range("U8").Value = 474
$ws.Range("W8").Value = 3.41
$ws.Range("X8").Value = 2.25
$ws.Range("Y8").Value = 10.31
$ws.Range("Z8").Value = 4.7
$ws.Range("AA8").Value = 112.7
$ws.Range("AC8").Value = 19083
$ws.Range("AD8").Value = 11
$ws.Range("AE8").Value = 192997
$ws.Range("AF8").Value = 1.09
$ws.Range("AG8").Value = 3167
$ws.Range("AH8").Value = 1.51
$ws.Range("AI8").Value = 16.59

$ws.Range("D9").Value = 34490
$ws.Range("E9").Value = 1199
$ws.Range("G9").Value = 1049
$ws.Range("H9").Value = 801
$ws.Range("I9").Value = 811
$ws.Range("K9").Value = 16746
$ws.Range("L9").Value = 8620
$ws.Range("M9").Value = 8128
$ws.Range("N9").Value = 8128
$ws.Range("P9").Value = 191
$ws.Range("Q9").Value = 1148
$ws.Range("R9").Value = -814
$ws.Range("S9").Value = -202
$ws.Range("T9").Value = 692
$ws.Range("U9").Value = 527
$ws.Range("W9").Value = 3.48
$ws.Range("X9").Value = 2.32
$ws.Range("Y9").Value = 10.41
$ws.Range("Z9").Value = 4.91
$ws.Range("AA9").Value = 106.04
$ws.Range("AC9").Value = 21004
$ws.Range("AD9").Value = 10
$ws.Range("AE9").Value = 210626
$ws.Range("AF9").Value = 1
$ws.Range("AG9").Value = 3250
$ws.Range("AH9").Value = 1.55
$ws.Range("AI9").Value = 15.47
